$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update "Periodo Mora" column (E16:E19) so the periods are listed in
# ascending order (2409, 2410, 2411, 2412) instead of descending.
$ws.Range("E16").Value = "2409"
$ws.Range("E17").Value = "2410"
$ws.Range("E18").Value = "2411"
$ws.Range("E19").Value = "2412"

# Update "Valor Mora" column (G16:G19) with the new amount.
$ws.Range("G16").Value = 2050000
$ws.Range("G17").Value = 2050000
$ws.Range("G18").Value = 2050000
$ws.Range("G19").Value = 2050000

$wb.Save()
